$d = $word.ActiveDocument

# --- Locate the run containing the lab-number digit ("4") ---
$find = $d.Content
$find.Find.ClearFormatting()
$null = $find.Find.Execute("4")
$digitStart = $find.Start
$digitEnd = $find.End

# --- Remove the existing "_GoBack" bookmark (it currently sits after
#     "Установите Honeyd:") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete() | Out-Null
}

# --- Change the lab number from 4 to 7, preserving the run/formatting ---
$digitRange = $d.Range($digitStart, $digitEnd)
$digitRange.Text = "7"

# --- Re-create the "_GoBack" bookmark right after the new "7", at the
#     end of that paragraph's text. A directly-collapsed Range cannot be
#     used here, so temporarily insert a marker character, wrap it with
#     the bookmark, then delete the marker -- the bookmark collapses back
#     to the correct position. ---
$marker = $d.Range($digitEnd, $digitEnd)
$marker.InsertAfter("~")

$wrapRange = $d.Range($digitEnd, $digitEnd + 1)
$d.Bookmarks.Add("_GoBack", $wrapRange)

$markerRange = $d.Range($digitEnd, $digitEnd + 1)
$markerRange.Text = ""
